$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: Array Rearrangement #3 (no hyperlink, plain URL text) ---
$ws.Range("C14").Value2 = "done"
$ws.Range("F14").Value2 = "Array Rearrangement #3"
$ws.Range("G14").Value2 = "https://www.geeksforgeeks.org/rearrange-array-arri-arrj-even-arri/"

# --- Row 15: Array Rearrangement #4 (with hyperlink) ---
$ws.Range("C15").Value2 = "done"
$ws.Range("F15").Value2 = "Array Rearrangement #4"
$ws.Range("G15").Value2 = "https://www.geeksforgeeks.org/rearrange-positive-and-negative-numbers-publish/"

$ws.Hyperlinks.Add($ws.Range("G15"), "https://www.geeksforgeeks.org/rearrange-positive-and-negative-numbers-publish/") | Out-Null
$ws.Range("G15").Style = "Hyperlink"

# --- Update the view: scroll/freeze pane + active selection ---
$ws.Range("C16").Select()
